$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "2025-09-26"
$ws.Range("A42").ClearFormats()

$ws.Range("B42").Value = "21:24:23"
$ws.Range("C42").Value = "1.00 EUR = 1,619.7750"
